$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") contains a date serial number for rows 2 through 206.
# Update all of these from 45202 (2023-10-03) to 45203 (2023-10-04).
$range = $ws.Range("C2:C206")
$range.Value = 45203
